$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "minutes_to_run_mile"
$ws.Range("B2").Select()
